$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" field text (1/7/2017 -> 10/25/18)
#    across the notes master, the slide master and every slide layout.
# ---------------------------------------------------------------------
$oldDate = "1/7/2017"
$newDate = "10/25/18"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Notes master date placeholder.
Update-DatePlaceholder $p.NotesMaster.Shapes

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's date placeholder.
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Rename the Person* shapes to Loan* on the diagram slide.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -eq -1) {
        $text = $sh.TextFrame.TextRange.Text
        if ($text -eq "PersonListPanel") {
            $sh.TextFrame.TextRange.Text = "LoanListPanel"
        } elseif ($text -eq "PersonCard") {
            $sh.TextFrame.TextRange.Text = "LoanCard"
            $sh.TextFrame.TextRange.Font.Bold = -1
        }
    }
}
